$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1554434735375247
$ws.Cells.Item(2, 3).Value = 0.3375848360084654
$ws.Cells.Item(2, 4).Value = 2938.103010863317
$ws.Cells.Item(2, 5).Value = 6.48142807727062
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2945.077467250134

$ws.Cells.Item(3, 2).Value = 0.7287194209349384
$ws.Cells.Item(3, 3).Value = 1.65323645889881
$ws.Cells.Item(3, 4).Value = 0.1529057820181812
$ws.Cells.Item(3, 5).Value = 6.48142807727062
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 9.016289739122548

$ws.Cells.Item(4, 2).Value = 1.505614041169197
$ws.Cells.Item(4, 3).Value = 9.226618575922256
$ws.Cells.Item(4, 4).Value = 16.98373111632243
$ws.Cells.Item(4, 5).Value = 6.48142807727062
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 34.1973918106845

$ws.Cells.Item(5, 2).Value = 0.7287194209349384
$ws.Cells.Item(5, 3).Value = 0.3375848360084654
$ws.Cells.Item(5, 4).Value = 3.082599426703578
$ws.Cells.Item(5, 5).Value = 0.4998867070740569
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.64879039072104

$ws.Cells.Item(6, 2).Value = 0.3464964993005633
$ws.Cells.Item(6, 3).Value = 1.65323645889881
$ws.Cells.Item(6, 4).Value = 0.7127328510149897
$ws.Cells.Item(6, 5).Value = 0.4998867070740569
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 3.21235251628842

$ws.Cells.Item(7, 2).Value = 1.505614041169197
$ws.Cells.Item(7, 3).Value = 1.65323645889881
$ws.Cells.Item(7, 4).Value = 0.7127328510149897
$ws.Cells.Item(7, 5).Value = 6.48142807727062
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 10.35301142835362

$ws.Cells.Item(8, 2).Value = 1.505614041169197
$ws.Cells.Item(8, 3).Value = 1.65323645889881
$ws.Cells.Item(8, 4).Value = 0.1529057820181812
$ws.Cells.Item(8, 5).Value = 0.4998867070740569
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.811642989160245

$ws.Cells.Item(9, 2).Value = 0.001754667048134761
$ws.Cells.Item(9, 3).Value = 0.05231270169004087
$ws.Cells.Item(9, 4).Value = 0.7127328510149897
$ws.Cells.Item(9, 5).Value = 0.4998867070740569
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.266686926827222

$ws.Cells.Item(10, 2).Value = 3.182878228561681
$ws.Cells.Item(10, 3).Value = 1.65323645889881
$ws.Cells.Item(10, 4).Value = 0.7127328510149897
$ws.Cells.Item(10, 5).Value = 6.48142807727062
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 12.0302756157461

$ws.Cells.Item(11, 2).Value = 0.7287194209349384
$ws.Cells.Item(11, 3).Value = 1.65323645889881
$ws.Cells.Item(11, 4).Value = 3.082599426703578
$ws.Cells.Item(11, 5).Value = 0.4998867070740569
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 5.964442013611383

